# Update factsheets with text edits from COMM
#
# The "No. of 990 Filers w/ Gov Grants" numeric columns across the
# workbook are converted from numbers to plain text (matching the
# formatting already used by every other column in these sheets), the
# overall filer count is written with a thousands separator ("1,561"),
# and a new "Total" row is added to the County sheet.

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to be treated as text so strings that look like
    # numbers (e.g. "420", "1,561", "$3,847,910,957", "7.14%") are not
    # silently re-interpreted as numeric/currency/percentage values.
    $range.NumberFormat = "@"
    $range.Value = $value
    # Drop back to the default style so no stray per-cell number format
    # is left behind once the text is safely stored.
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overall": A2 1561 -> "1,561"
# ---------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall.Range("A2") "1,561"

# ---------------------------------------------------------------
# Sheet "County": B2:B10 numbers -> text, plus new Total row 11
# ---------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")
Set-TextValue $wsCounty.Range("B2") "420"
Set-TextValue $wsCounty.Range("B3") "109"
Set-TextValue $wsCounty.Range("B4") "90"
Set-TextValue $wsCounty.Range("B5") "111"
Set-TextValue $wsCounty.Range("B6") "21"
Set-TextValue $wsCounty.Range("B7") "91"
Set-TextValue $wsCounty.Range("B8") "280"
Set-TextValue $wsCounty.Range("B9") "135"
Set-TextValue $wsCounty.Range("B10") "304"

Set-TextValue $wsCounty.Range("A11") "Total"
Set-TextValue $wsCounty.Range("B11") "1,561"
Set-TextValue $wsCounty.Range("C11") "$3,847,910,957"
Set-TextValue $wsCounty.Range("D11") "7.14%"
Set-TextValue $wsCounty.Range("E11") "-16.67%"
Set-TextValue $wsCounty.Range("F11") "71.11%"

# ---------------------------------------------------------------
# Sheet "Congressional District": B2:B6 numbers -> text, B7 Total -> "1,561"
# ---------------------------------------------------------------
$wsCd = $wb.Worksheets.Item("Congressional District")
Set-TextValue $wsCd.Range("B2") "359"
Set-TextValue $wsCd.Range("B3") "264"
Set-TextValue $wsCd.Range("B4") "304"
Set-TextValue $wsCd.Range("B5") "339"
Set-TextValue $wsCd.Range("B6") "295"
Set-TextValue $wsCd.Range("B7") "1,561"

# ---------------------------------------------------------------
# Sheet "Size": B2:B7 numbers -> text, B8 Total -> "1,561"
# ---------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
Set-TextValue $wsSize.Range("B2") "436"
Set-TextValue $wsSize.Range("B3") "446"
Set-TextValue $wsSize.Range("B4") "274"
Set-TextValue $wsSize.Range("B5") "124"
Set-TextValue $wsSize.Range("B6") "195"
Set-TextValue $wsSize.Range("B7") "86"
Set-TextValue $wsSize.Range("B8") "1,561"

# ---------------------------------------------------------------
# Sheet "Subsector": B2:B13 numbers -> text, B14 Total -> "1,561"
# ---------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")
Set-TextValue $wsSub.Range("B2") "176"
Set-TextValue $wsSub.Range("B3") "206"
Set-TextValue $wsSub.Range("B4") "53"
Set-TextValue $wsSub.Range("B5") "149"
Set-TextValue $wsSub.Range("B6") "1"
Set-TextValue $wsSub.Range("B7") "456"
Set-TextValue $wsSub.Range("B8") "14"
Set-TextValue $wsSub.Range("B9") "1"
Set-TextValue $wsSub.Range("B10") "110"
Set-TextValue $wsSub.Range("B11") "16"
Set-TextValue $wsSub.Range("B12") "365"
Set-TextValue $wsSub.Range("B13") "14"
Set-TextValue $wsSub.Range("B14") "1,561"

Write-Output "edits applied"
